$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add today's status rows (22-25) for 18-12-2024 ---

# Copy the date-cell formatting (numeric date format + center/center alignment,
# same style already used by A3/A6/A10/A14/A18) down onto the new date cell A22.
$ws.Range("A21").Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A22").Value = 45644                      # 18-Dec-2024

$ws.Range("B22").Value = "ppt documents(3 sheets) ,2 incomplete"
$ws.Range("E22").Value = "Completed"

$ws.Range("B23").Value = "sheet preparation for leave tracker "
$ws.Range("E23").Value = "Completed"

$ws.Range("B24").Value = "login late tracker generation in tableau(sample file generation)"
$ws.Range("E24").Value = "Completed"

$ws.Range("B25").Value = "sheet preparation for leave tracker "
$ws.Range("E25").Value = "Completed"

# Merge the date column for the new day's block, matching the existing
# A3:A5 / A6:A9 / A10:A13 / A14:A17 / A18:A21 pattern.
$ws.Range("A22:A25").Merge() | Out-Null

# --- Update the cursor/selection to reflect where the user ended up ---
$ws.Range("B20").Select() | Out-Null
